# Update workbook for "Add data for 2022-04-25" commit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / tab title to reflect new "through" date
$ws.Name = "Through 2022-04-17"

# Update the label cell for April row to reflect the new "through" date
$ws.Range("A5").Value = "April (through 04-17)"

# Update April row (row 5) values for years 2015-2022 (columns B-I)
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 34
$ws.Range("E5").Value = 29
$ws.Range("F5").Value = 28
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = 63
$ws.Range("I5").Value = 74

# Update Total row (row 6) values for years 2015-2022 (columns B-I)
$ws.Range("B6").Value = 78
$ws.Range("C6").Value = 142
$ws.Range("D6").Value = 223
$ws.Range("E6").Value = 226
$ws.Range("F6").Value = 138
$ws.Range("G6").Value = 236
$ws.Range("H6").Value = 486
$ws.Range("I6").Value = 509
